# [PROS-1480] Removal of Hyphens in scene type
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace hyphens with spaces in the scene-type / store_area_location strings (column D, rows 2-10)
$ws.Range("D2").Value = "Event Space, Checkout, Beverage end, Deli section, Liquor, 0003 Event space, 0004 Checkout, 0005 End, 0006 Deli, 0007 Liquor"
$ws.Range("D3").Value = "Event Space, Checkout, Deli section, Liquor, 0003 Event space,0004 Checkout, 0006 Deli, 0007 Liquor"
$ws.Range("D4").Value = "Event Space, Checkout, Deli section, 0003 Event space, 0004 Checkout, 0006 Deli"
$ws.Range("D5").Value = "Event Space, Checkout, Beverage end, Deli section, 0003 Event space,0004 Checkout, 0005 End, 0006 Deli"
$ws.Range("D6").Value = "Event Space, Checkout, Store front, 0003 Event space, 0004 Checkout, 0009 Store Front"
$ws.Range("D7").Value = "Checkout, Store front, 0004 Checkout, 0009 Store Front"
$ws.Range("D8").Value = "Event Space, Checkout, Beverage end, Liquor, 0003 Event space, 0004 Checkout, 0005 End, 0007 Liquor"
$ws.Range("D9").Value = "Event Space, Checkout, Beverage end, 0003 Event space, 0004 Checkout, 0005 End"
$ws.Range("D10").Value = "Event Space, Checkout, 0003 Event space, 0004 Checkout"

# Update the view: zoom out to 90% and move the selection to D10 (as in the saved view state)
$win = $ws.Application.ActiveWindow
$win.Zoom = 90
$ws.Range("D10").Select()
